$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: new columns D (email, hyperlinked) and E (note) ---
$ws.Range("D1").Value = "espoirditekemena@gmail.com"
$ws.Hyperlinks.Add($ws.Range("D1"), "mailto:espoirditekemena@gmail.com")
$ws.Range("E1").Value = "jds;oaiudisuygfoiu"

# --- Row 2: new record (Miriam wa mbuyi) ---
$ws.Range("A2").Value = "Miriam"
$ws.Range("B2").Value = "Mbuyi"
$ws.Range("C2").Value = "wa mbuyi"
$ws.Range("D2").Value = "miriamMbuyi@gmail.com"
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:miriamMbuyi@gmail.com")
$ws.Range("E2").Value = "jdhfdsfhjhffjdhdsj"

# --- Column widths for the new columns ---
$ws.Columns.Item(4).ColumnWidth = 32.6
$ws.Columns.Item(5).ColumnWidth = 24.6

# --- Selection matches the authored state ---
$ws.Range("E12").Select()
